# Auto-generated Excel COM-interop edit script
# Updates horarios-141-2026-01-11.xlsx per the scraper refresh diff

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LP1912")

# --- Header: last-updated timestamp + row count (sheet "LP1912") ---
$ws.Cells.Item(2,1).Value = "Última actualización: 16:44:12"
$ws.Cells.Item(3,1).Value = "Total filas: 267"

# --- Rows 106-108: re-sorted (tie on Hora_Llegada 11:52) ---
$ws.Cells.Item(106,1).Value = "10:05:51"
$ws.Cells.Item(106,2).Value = "11:52"
$ws.Cells.Item(106,3).Value = "225_GOMEZ"
$ws.Cells.Item(106,4).Value = 107
$ws.Cells.Item(106,5).Value = "LP1912"
$ws.Cells.Item(107,1).Value = "11:47:17"
$ws.Cells.Item(107,2).Value = "11:52"
$ws.Cells.Item(107,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(107,4).Value = 5
$ws.Cells.Item(107,5).Value = "LP1912"
$ws.Cells.Item(108,1).Value = "11:52:01"
$ws.Cells.Item(108,2).Value = "11:52"
$ws.Cells.Item(108,3).Value = "15X38_ABASTO"
$ws.Cells.Item(108,4).Value = 0
$ws.Cells.Item(108,5).Value = "LP1912"

# --- Rows 137-138: re-sorted (tie on Hora_Llegada 12:37) ---
$ws.Cells.Item(137,1).Value = "11:47:17"
$ws.Cells.Item(137,2).Value = "12:37"
$ws.Cells.Item(137,3).Value = "27_EL RETIRO"
$ws.Cells.Item(137,4).Value = 50
$ws.Cells.Item(137,5).Value = "LP1912"
$ws.Cells.Item(138,1).Value = "11:52:01"
$ws.Cells.Item(138,2).Value = "12:37"
$ws.Cells.Item(138,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(138,4).Value = 45
$ws.Cells.Item(138,5).Value = "LP1912"

# --- Rows 139-141: re-sorted (tie on Hora_Llegada 12:47, Minutos 73) ---
$ws.Cells.Item(139,1).Value = "11:34:59"
$ws.Cells.Item(139,2).Value = "12:47"
$ws.Cells.Item(139,3).Value = "14_ABASTO"
$ws.Cells.Item(139,4).Value = 73
$ws.Cells.Item(139,5).Value = "LP1912"
$ws.Cells.Item(140,1).Value = "11:34:59"
$ws.Cells.Item(140,2).Value = "12:47"
$ws.Cells.Item(140,3).Value = "15X38_ABASTO"
$ws.Cells.Item(140,4).Value = 73
$ws.Cells.Item(140,5).Value = "LP1912"
$ws.Cells.Item(141,1).Value = "11:34:59"
$ws.Cells.Item(141,2).Value = "12:47"
$ws.Cells.Item(141,3).Value = "16_SANTA ANA"
$ws.Cells.Item(141,4).Value = 73
$ws.Cells.Item(141,5).Value = "LP1912"

# --- New scrape run added two more arrivals near the end of the day;
#     insert two blank rows at 270 (pushes old 270 -> 272), then
#     rewrite rows 242..272 with the refreshed, re-sorted data ---
$ws.Rows.Item(270).Insert()
$ws.Rows.Item(270).Insert()

$ws.Cells.Item(242,1).Value = "16:44:12"
$ws.Cells.Item(242,2).Value = "17:08"
$ws.Cells.Item(242,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(242,4).Value = 24
$ws.Cells.Item(242,5).Value = "LP1912"

$ws.Cells.Item(244,1).Value = "15:46:07"
$ws.Cells.Item(244,2).Value = "17:09"
$ws.Cells.Item(244,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(244,4).Value = 83
$ws.Cells.Item(244,5).Value = "LP1912"

$ws.Cells.Item(245,1).Value = "15:17:33"
$ws.Cells.Item(245,2).Value = "17:10"
$ws.Cells.Item(245,3).Value = "215C_EL PATO"
$ws.Cells.Item(245,4).Value = 113
$ws.Cells.Item(245,5).Value = "LP1912"

$ws.Cells.Item(246,1).Value = "15:46:07"
$ws.Cells.Item(246,2).Value = "17:21"
$ws.Cells.Item(246,3).Value = "15X38_ABASTO"
$ws.Cells.Item(246,4).Value = 95
$ws.Cells.Item(246,5).Value = "LP1912"

$ws.Cells.Item(247,1).Value = "16:28:03"
$ws.Cells.Item(247,2).Value = "17:28"
$ws.Cells.Item(247,3).Value = "27_EL RETIRO"
$ws.Cells.Item(247,4).Value = 60
$ws.Cells.Item(247,5).Value = "LP1912"

$ws.Cells.Item(248,1).Value = "16:14:44"
$ws.Cells.Item(248,2).Value = "17:32"
$ws.Cells.Item(248,3).Value = "27_EL RETIRO"
$ws.Cells.Item(248,4).Value = 78
$ws.Cells.Item(248,5).Value = "LP1912"

$ws.Cells.Item(249,1).Value = "15:46:07"
$ws.Cells.Item(249,2).Value = "17:34"
$ws.Cells.Item(249,3).Value = "17_ROMERO"
$ws.Cells.Item(249,4).Value = 108
$ws.Cells.Item(249,5).Value = "LP1912"

$ws.Cells.Item(250,1).Value = "15:58:05"
$ws.Cells.Item(250,2).Value = "17:36"
$ws.Cells.Item(250,3).Value = "27_EL RETIRO"
$ws.Cells.Item(250,4).Value = 98
$ws.Cells.Item(250,5).Value = "LP1912"

$ws.Cells.Item(251,1).Value = "15:46:07"
$ws.Cells.Item(251,2).Value = "17:37"
$ws.Cells.Item(251,3).Value = "27_EL RETIRO"
$ws.Cells.Item(251,4).Value = 111
$ws.Cells.Item(251,5).Value = "LP1912"

$ws.Cells.Item(252,1).Value = "16:14:44"
$ws.Cells.Item(252,2).Value = "17:38"
$ws.Cells.Item(252,3).Value = "215B_EL PATO"
$ws.Cells.Item(252,4).Value = 84
$ws.Cells.Item(252,5).Value = "LP1912"

$ws.Cells.Item(253,1).Value = "15:46:07"
$ws.Cells.Item(253,2).Value = "17:39"
$ws.Cells.Item(253,3).Value = "215B_EL PATO"
$ws.Cells.Item(253,4).Value = 113
$ws.Cells.Item(253,5).Value = "LP1912"

$ws.Cells.Item(254,1).Value = "16:14:44"
$ws.Cells.Item(254,2).Value = "17:42"
$ws.Cells.Item(254,3).Value = "215_EL PELIGRO"
$ws.Cells.Item(254,4).Value = 116
$ws.Cells.Item(254,5).Value = "LP1912"

$ws.Cells.Item(255,1).Value = "16:14:44"
$ws.Cells.Item(255,2).Value = "17:45"
$ws.Cells.Item(255,3).Value = "215_EL PELIGRO"
$ws.Cells.Item(255,4).Value = 91
$ws.Cells.Item(255,5).Value = "LP1912"

$ws.Cells.Item(256,1).Value = "15:58:05"
$ws.Cells.Item(256,2).Value = "17:46"
$ws.Cells.Item(256,3).Value = "215_EL PELIGRO"
$ws.Cells.Item(256,4).Value = 108
$ws.Cells.Item(256,5).Value = "LP1912"

$ws.Cells.Item(257,1).Value = "16:44:12"
$ws.Cells.Item(257,2).Value = "17:48"
$ws.Cells.Item(257,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(257,4).Value = 64
$ws.Cells.Item(257,5).Value = "LP1912"

$ws.Cells.Item(258,1).Value = "16:28:03"
$ws.Cells.Item(258,2).Value = "17:49"
$ws.Cells.Item(258,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(258,4).Value = 81
$ws.Cells.Item(258,5).Value = "LP1912"

$ws.Cells.Item(259,1).Value = "16:37:06"
$ws.Cells.Item(259,2).Value = "17:50"
$ws.Cells.Item(259,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(259,4).Value = 73
$ws.Cells.Item(259,5).Value = "LP1912"

$ws.Cells.Item(260,1).Value = "16:37:06"
$ws.Cells.Item(260,2).Value = "17:52"
$ws.Cells.Item(260,3).Value = "10_OLMOS"
$ws.Cells.Item(260,4).Value = 75
$ws.Cells.Item(260,5).Value = "LP1912"

$ws.Cells.Item(261,1).Value = "16:28:03"
$ws.Cells.Item(261,2).Value = "17:53"
$ws.Cells.Item(261,3).Value = "10_OLMOS"
$ws.Cells.Item(261,4).Value = 85
$ws.Cells.Item(261,5).Value = "LP1912"

$ws.Cells.Item(262,1).Value = "16:28:03"
$ws.Cells.Item(262,2).Value = "17:58"
$ws.Cells.Item(262,3).Value = "17_ROMERO"
$ws.Cells.Item(262,4).Value = 90
$ws.Cells.Item(262,5).Value = "LP1912"

$ws.Cells.Item(263,1).Value = "16:14:44"
$ws.Cells.Item(263,2).Value = "18:00"
$ws.Cells.Item(263,3).Value = "10_OLMOS"
$ws.Cells.Item(263,4).Value = 106
$ws.Cells.Item(263,5).Value = "LP1912"

$ws.Cells.Item(264,1).Value = "16:28:03"
$ws.Cells.Item(264,2).Value = "18:05"
$ws.Cells.Item(264,3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(264,4).Value = 111
$ws.Cells.Item(264,5).Value = "LP1912"

$ws.Cells.Item(265,1).Value = "16:28:03"
$ws.Cells.Item(265,2).Value = "18:06"
$ws.Cells.Item(265,3).Value = "15_ABASTO"
$ws.Cells.Item(265,4).Value = 98
$ws.Cells.Item(265,5).Value = "LP1912"

$ws.Cells.Item(266,1).Value = "16:14:44"
$ws.Cells.Item(266,2).Value = "18:10"
$ws.Cells.Item(266,3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(266,4).Value = 116
$ws.Cells.Item(266,5).Value = "LP1912"

$ws.Cells.Item(267,1).Value = "16:14:44"
$ws.Cells.Item(267,2).Value = "18:10"
$ws.Cells.Item(267,3).Value = "15_ABASTO"
$ws.Cells.Item(267,4).Value = 116
$ws.Cells.Item(267,5).Value = "LP1912"

$ws.Cells.Item(268,1).Value = "16:28:03"
$ws.Cells.Item(268,2).Value = "18:17"
$ws.Cells.Item(268,3).Value = "10_OLMOS"
$ws.Cells.Item(268,4).Value = 109
$ws.Cells.Item(268,5).Value = "LP1912"

$ws.Cells.Item(269,1).Value = "16:28:03"
$ws.Cells.Item(269,2).Value = "18:22"
$ws.Cells.Item(269,3).Value = "215C_EL PATO"
$ws.Cells.Item(269,4).Value = 114
$ws.Cells.Item(269,5).Value = "LP1912"

$ws.Cells.Item(270,1).Value = "16:28:03"
$ws.Cells.Item(270,2).Value = "18:25"
$ws.Cells.Item(270,3).Value = "16_SANTA ANA"
$ws.Cells.Item(270,4).Value = 117
$ws.Cells.Item(270,5).Value = "LP1912"

$ws.Cells.Item(271,1).Value = "16:37:06"
$ws.Cells.Item(271,2).Value = "18:30"
$ws.Cells.Item(271,3).Value = "14_ABASTO"
$ws.Cells.Item(271,4).Value = 113
$ws.Cells.Item(271,5).Value = "LP1912"

$ws.Cells.Item(272,1).Value = "16:37:06"
$ws.Cells.Item(272,2).Value = "18:36"
$ws.Cells.Item(272,3).Value = "15X38_ABASTO"
$ws.Cells.Item(272,4).Value = 119
$ws.Cells.Item(272,5).Value = "LP1912"

# --- Propagate the refreshed timestamp to the other two sheets ---
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Cells.Item(2,1).Value = "Última actualización: 16:44:12"

$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Cells.Item(2,1).Value = "Última actualización: 16:44:12"
